# Applies the "updating book and updated syllabus" edit to the
# GEOG 473/673 spring syllabus.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $old"
    }
}

# ---------------------------------------------------------------------
# 1. Turn the "simkins@udel.edu" mailto field into a real w:hyperlink
#    (matching the already-hyperlinked "Website:" link next to it).
# ---------------------------------------------------------------------
$mailField = $null
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $f = $d.Fields.Item($i)
    if ($f.Code.Text -match "mailto:simkins@udel.edu") {
        $mailField = $f
        break
    }
}
if ($mailField -ne $null) {
    # Unlink() removes the field code but keeps the displayed result text
    # ("simkins@udel.edu") in place, with its existing Hyperlink character style.
    $mailField.Unlink()
}

$rng = $d.Content
$rng.Find.Execute("simkins@udel.edu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($rng, "mailto:simkins@udel.edu", "", "", "simkins@udel.edu") | Out-Null

# ---------------------------------------------------------------------
# 2. Coding Assignments paragraph: Python section -> Advanced R 2-credit
#    section, 33% -> 20% "of this portions grade."
# ---------------------------------------------------------------------
Replace-Text "For the Python section, a final project will be given which represents 33% of this portion." "For the Advanced R 2-credit section, a final project will be given which represents 20% of this portions grade."

# ---------------------------------------------------------------------
# 3. Remove the "Python Programming and Visualization for Scientists by
#    Alex J. Decaria" bullet from the textbook list.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Decaria") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 4. Course Outline - Phase 1 schedule updates
# ---------------------------------------------------------------------
Replace-Text "2/10/2020 to 4/17/2020" "2/15/2021 to 4/16/2021"

Replace-Text "2/15: Welcome to advanced R - Tutorial" "2/15: Plotting Customization Tutorial"
Replace-Text "2/17: Welcome to advanced R - Tutorial" "2/17: Plotting Customization Tutorial"

Replace-Text "4/5: R final project" "4/5: Advanced R project"
Replace-Text "4/7: R final project" "4/7: Advanced R project"
Replace-Text "4/12: R final project" "4/12: Advanced R project"
Replace-Text "4/14: R final project" "4/14: Advanced R project"

# ---------------------------------------------------------------------
# 5. Course Outline - Phase 2 header dates
# ---------------------------------------------------------------------
Replace-Text "Investigative Spatial Programming - 4/20/2020 to 5/18/2020  " "Investigative Spatial Programming - 4/19/2021 to 5/18/2021"

# ---------------------------------------------------------------------
# 6. Course Outline - Phase 2 schedule updates
# ---------------------------------------------------------------------
Replace-Text "4/19: Introduction to Investigative Spatial Programming – A group discussion" "4/19: Introduction to Machine Learning"
Replace-Text "4/21: Create Outline" "4/21: Introduction to Machine Learning"
Replace-Text "4/26: Project" "4/26: Time Series Forecasting "
Replace-Text "4/28: Project" "4/28: Time Series Forecasting"
Replace-Text "5/3: Project" "5/3: Time Series Forecasting"
Replace-Text "5/5: Project" "5/5: Time Series Forecasting / Random Forest Modeling"
Replace-Text "5/10: Project" "5/10: Random Forest Modeling"
Replace-Text "5/12: Project & Final Evaluation" "5/12: Random Forest Modeling"

# Add a brand-new trailing paragraph: "5/17: Random Forest Modeling"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^5/12: Random Forest Modeling") {
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = "5/17: Random Forest Modeling"
        break
    }
}

Write-Output "Edit complete."
